$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1993464052287582
$ws.Range("C2").Value = 0.565359477124183
$ws.Range("J2").Value = 0.0196078431372549
$ws.Range("P2").Value = 0.1339869281045752
$ws.Range("S2").Value = 0.08169934640522876
$ws.Range("B3").Value = 0.01092896174863388
$ws.Range("C3").Value = 0.0273224043715847
$ws.Range("J3").Value = 0.03278688524590164
$ws.Range("P3").Value = 0.7650273224043715
$ws.Range("S3").Value = 0.1639344262295082
$ws.Range("J4").Value = 0.04545454545454546
$ws.Range("P4").Value = 0.7045454545454546
$ws.Range("B6").Value = 0.08374384236453201
$ws.Range("F6").Value = 0.05911330049261083
$ws.Range("J6").Value = 0.2758620689655172
$ws.Range("O6").Value = 0.009852216748768473
$ws.Range("Q6").Value = 0.1576354679802956
$ws.Range("R6").Value = 0.07389162561576355
$ws.Range("S6").Value = 0.3399014778325123
$ws.Range("B7").Value = 0.1194029850746269
$ws.Range("D7").Value = 0.02238805970149254
$ws.Range("F7").Value = 0.05970149253731343
$ws.Range("J7").Value = 0.09701492537313433
$ws.Range("O7").Value = 0.02238805970149254
$ws.Range("Q7").Value = 0.2462686567164179
$ws.Range("R7").Value = 0.06716417910447761
$ws.Range("S7").Value = 0.3656716417910448
$ws.Range("B8").Value = 0.1057692307692308
$ws.Range("D8").Value = 0.01682692307692308
$ws.Range("F8").Value = 0.05288461538461538
$ws.Range("J8").Value = 0.1105769230769231
$ws.Range("O8").Value = 0.01682692307692308
$ws.Range("Q8").Value = 0.1995192307692308
$ws.Range("R8").Value = 0.1129807692307692
$ws.Range("S8").Value = 0.3846153846153846
$ws.Range("B9").Value = 0.06936416184971098
$ws.Range("D9").Value = 0.02312138728323699
$ws.Range("F9").Value = 0.06358381502890173
$ws.Range("J9").Value = 0.1329479768786127
$ws.Range("O9").Value = 0.01734104046242774
$ws.Range("Q9").Value = 0.1676300578034682
$ws.Range("R9").Value = 0.115606936416185
$ws.Range("S9").Value = 0.4104046242774567
$ws.Range("B10").Value = 0.1222596964586847
$ws.Range("D10").Value = 0.02529510961214165
$ws.Range("F10").Value = 0.06155143338954469
$ws.Range("J10").Value = 0.1104553119730186
$ws.Range("O10").Value = 0.01686340640809443
$ws.Range("Q10").Value = 0.2445193929173693
$ws.Range("R10").Value = 0.1037099494097808
$ws.Range("S10").Value = 0.315345699831366
$ws.Range("G11").Value = 0.1271929824561404
$ws.Range("J11").Value = 0.09649122807017543
$ws.Range("K11").Value = 0.1929824561403509
$ws.Range("L11").Value = 0.5657894736842105
$ws.Range("S11").Value = 0.01754385964912281
$ws.Range("G12").Value = 0.7196969696969697
$ws.Range("J12").Value = 0.1893939393939394
$ws.Range("K12").Value = 0.007575757575757576
$ws.Range("L12").Value = 0.03787878787878788
$ws.Range("S12").Value = 0.04545454545454546
$ws.Range("G13").Value = 0.46875
$ws.Range("J13").Value = 0.40625
$ws.Range("S13").Value = 0.125
$ws.Range("F15").Value = 0.04225352112676056
$ws.Range("H15").Value = 0.1596244131455399
$ws.Range("I15").Value = 0.06103286384976526
$ws.Range("J15").Value = 0.3849765258215962
$ws.Range("K15").Value = 0.05164319248826291
$ws.Range("M15").Value = 0.004694835680751174
$ws.Range("O15").Value = 0.07042253521126761
$ws.Range("S15").Value = 0.2253521126760563
$ws.Range("F16").Value = 0.01477832512315271
$ws.Range("H16").Value = 0.2068965517241379
$ws.Range("I16").Value = 0.07389162561576355
$ws.Range("J16").Value = 0.354679802955665
$ws.Range("K16").Value = 0.09359605911330049
$ws.Range("M16").Value = 0.009852216748768473
$ws.Range("O16").Value = 0.09359605911330049
$ws.Range("S16").Value = 0.1527093596059113
$ws.Range("F17").Value = 0.02132196162046908
$ws.Range("H17").Value = 0.1940298507462687
$ws.Range("I17").Value = 0.1108742004264392
$ws.Range("J17").Value = 0.4008528784648188
$ws.Range("K17").Value = 0.08315565031982942
$ws.Range("M17").Value = 0.01492537313432836
$ws.Range("O17").Value = 0.06609808102345416
$ws.Range("S17").Value = 0.1087420042643923
$ws.Range("F18").Value = 0.009523809523809525
$ws.Range("H18").Value = 0.2476190476190476
$ws.Range("I18").Value = 0.04285714285714286
$ws.Range("J18").Value = 0.4476190476190476
$ws.Range("K18").Value = 0.05714285714285714
$ws.Range("M18").Value = 0.01428571428571429
$ws.Range("O18").Value = 0.0761904761904762
$ws.Range("S18").Value = 0.1047619047619048
$ws.Range("F19").Value = 0.023963133640553
$ws.Range("H19").Value = 0.1824884792626728
$ws.Range("I19").Value = 0.07926267281105991
$ws.Range("J19").Value = 0.3870967741935484
$ws.Range("K19").Value = 0.09400921658986175
$ws.Range("M19").Value = 0.0184331797235023
$ws.Range("O19").Value = 0.06912442396313365
$ws.Range("S19").Value = 0.1456221198156682
